# UserStories.xlsx — "Updated Car Inventory Project."
#
# Content change: row 12's "As a(n) <actor>" cell combined the
# "Sales Representative" role with "Manager" into a single cell,
# "Sales Representative, Manager" (the old standalone "Sales
# Representative" shared string is no longer referenced anywhere
# once this is applied).
#
# View change: the last on-screen selection before save moved to B12
# (the description cell of that same row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the actor cell for the "Manage the car inventory..." user story.
$ws.Range("A12").Value = "Sales Representative, Manager"

# Reflect the workbook's final selection/cursor position.
$ws.Range("B12").Select()
